# This script applies a weekly update to the "Hortaliza, Feria Lagunitas de
# Puerto Montt - Tomate" sheet: three brand-new price rows are inserted into
# the existing data table (which pushes the rows below them down), while all
# pre-existing rows keep their original values.
#
# Layout recap (columns A..R):
# A Mercado ID | B Mercado | C Región | D Fecha | E Codreg | F Categoría ID
# G Categoría | H Variedad | I Calidad | J Volumen | K Precio minimo
# L Precio maximo | M Precio promedio ponderado | N Unidad de comercializacion
# O Origen | P Precio $/Kg | Q Kg o Unidades | R Clasificacion

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow {
    param(
        [int]$Row,
        [double]$Fecha,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [string]$Unidad,
        [string]$Origen,
        [double]$PrecioKg,
        [double]$KgUnidades
    )

    $ws.Range("A$Row").Value = 4
    $ws.Range("B$Row").Value = "Feria Lagunitas de Puerto Montt"
    $ws.Range("C$Row").Value = "Los Lagos"
    $ws.Range("D$Row").Value = $Fecha
    $ws.Range("E$Row").Value = 10
    $ws.Range("F$Row").Value = 100112020
    $ws.Range("G$Row").Value = "Tomate"
    $ws.Range("H$Row").Value = "Larga vida"
    $ws.Range("I$Row").Value = $Calidad
    $ws.Range("J$Row").Value = $Volumen
    $ws.Range("K$Row").Value = $PrecioMin
    $ws.Range("L$Row").Value = $PrecioMax
    $ws.Range("M$Row").Value = $PrecioProm
    $ws.Range("N$Row").Value = $Unidad
    $ws.Range("O$Row").Value = $Origen
    $ws.Range("P$Row").Value = $PrecioKg
    $ws.Range("Q$Row").Value = $KgUnidades
    $ws.Range("R$Row").Value = "Hortaliza"
}

# --- Insert new row at 595 -------------------------------------------------
# All rows from 595 downward shift down by one.
$ws.Rows.Item(595).Insert()
Set-DataRow -Row 595 -Fecha 44748 -Calidad "Primera" -Volumen 120 `
    -PrecioMin 17000 -PrecioMax 17000 -PrecioProm 17000 `
    -Unidad "$/bandeja 18 kilos" -Origen "Región de Arica y Parinacota" `
    -PrecioKg 944 -KgUnidades 18

# --- Insert two new rows at 646-647 ----------------------------------------
# (Positions are expressed after the first insertion above.) Everything from
# row 646 downward shifts down by two more rows.
$ws.Range("A646:A647").EntireRow.Insert()
Set-DataRow -Row 646 -Fecha 44747 -Calidad "Primera" -Volumen 500 `
    -PrecioMin 17000 -PrecioMax 17000 -PrecioProm 17000 `
    -Unidad "$/bandeja 18 kilos" -Origen "Región de Arica y Parinacota" `
    -PrecioKg 944 -KgUnidades 18
Set-DataRow -Row 647 -Fecha 44747 -Calidad "Segunda" -Volumen 500 `
    -PrecioMin 13000 -PrecioMax 13000 -PrecioProm 13000 `
    -Unidad "$/bandeja 18 kilos" -Origen "Región de Arica y Parinacota" `
    -PrecioKg 722 -KgUnidades 18

$u = $ws.UsedRange.Address()
Write-Host "Final UsedRange: $u"
